$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (DAMSLTag, DialogAct)
$updates = @{
    18 = @("sd", "Statement-non-opinion")
    25 = @("sd", "Statement-non-opinion")
    32 = @("sd", "Statement-non-opinion")
    48 = @("aa", "Agree/Accept")
    70 = @("sv", "Statement-opinion")
    78 = @("sd", "Statement-non-opinion")
    83 = @("sv", "Statement-opinion")
    93 = @("sv", "Statement-opinion")
    97 = @("sv", "Statement-opinion")
    99 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
